$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits at the end
#    of the "...is not an integer)." paragraph (right after the period).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Relocate the "301 (Unauthorized)" heading + its body paragraph so
#    they land after the "400 (Bad request)" section (just before the
#    "404 (Not found)" heading), and rename the heading to
#    "401 (Unauthorized)".
# ------------------------------------------------------------------
$d.Content.Find.Execute("301 (Unauthorized)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "401 (Unauthorized)", 2) | Out-Null

$headingPara = $null
$bodyPara = $null
$notFoundPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq "401 (Unauthorized)`r") {
        $headingPara = $i
        $bodyPara = $i + 1
    }
    if ($txt -eq "404 (Not found)`r") {
        $notFoundPara = $i
    }
}

$moveStart = $d.Paragraphs.Item($headingPara).Range.Start
$moveEnd = $d.Paragraphs.Item($bodyPara).Range.End
$moveRange = $d.Range($moveStart, $moveEnd)
$moveRange.Cut()

# The cut removed two whole paragraphs that lived *before* the "404 (Not
# found)" heading, so its paragraph index shifted down by 2.
$notFoundPara = $notFoundPara - 2

$targetStart = $d.Paragraphs.Item($notFoundPara).Range.Start
$destRange = $d.Range($targetStart, $targetStart)
$destRange.Paste()

# Restore the Heading 3 style on the pasted heading paragraph (paste does
# not reliably keep the paragraph-mark formatting of a moved paragraph).
$d.Paragraphs.Item($notFoundPara).Style = "Heading 3"

# ------------------------------------------------------------------
# 3. Re-merge the body paragraph's three runs ("The credentials were " +
#    "invalid" + " or not present.") into a single run of plain text -
#    simplest way is just to rewrite the paragraph's text.
# ------------------------------------------------------------------
$bodyIndex = $notFoundPara + 1
$d.Paragraphs.Item($bodyIndex).Range.Text = "The credentials were invalid or not present.`r"

# ------------------------------------------------------------------
# 4. Re-add the "_GoBack" bookmark inside the "400 (Bad request)" heading
#    text, splitting it into "400 (Bad requ" | bookmark | "est)" -
#    mirroring where the cursor was left after the last edit.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq "400 (Bad request)`r") {
        $pos = $d.Paragraphs.Item($i).Range.Start + 13
        $bmRange = $d.Range($pos, $pos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
